# Refresh cryptocurrency price / volume(1h) data (and a couple of
# re-ordered / replaced coin rows) as scraped by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 'Price' column (D) holds numbers formatted with '.' as a thousands
# separator as plain text (e.g. "27.434.96"), so force text formatting
# first to stop Excel from re-interpreting values such as "332.96" as a
# real number, then drop back to the default style once written so no
# stray cell formatting is introduced.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.434.96'
$ws.Range('E2').Value = '  -2.10%  '
$ws.Range('D3').Value = '1.835.68'
$ws.Range('E3').Value = '  -2.62%  '
$ws.Range('E4').Value = '  -0.99%  '
$ws.Range('D5').Value = '332.96'
$ws.Range('E5').Value = '  -0.75%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  -0.87%  '
$ws.Range('D7').Value = '0.4621'
$ws.Range('E7').Value = '  -2.57%  '
$ws.Range('D8').Value = '0.3815'
$ws.Range('E8').Value = '  -3.32%  '
$ws.Range('D9').Value = '46.53'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').Value = '0.07888'
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('D11').Value = '0.9725'
$ws.Range('E11').Value = '  -4.56%  '
$ws.Range('D12').Value = '21.08'
$ws.Range('E12').Value = '  -3.58%  '
$ws.Range('D13').Value = '1.848.24'
$ws.Range('E13').Value = '  -3.89%  '
$ws.Range('D14').Value = '5.902'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '7.032'
$ws.Range('E15').Value = '  -2.34%  '
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('D17').Value = '87.84'
$ws.Range('E17').Value = '  -0.75%  '
$ws.Range('D18').Value = '0.06616'
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('D19').Value = '0.00001027'
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('D20').Value = '16.95'
$ws.Range('E20').Value = '  -0.85%  '
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('D22').Value = '27.444.57'
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('D23').Value = '5.353'
$ws.Range('E23').Value = '  -2.82%  '
$ws.Range('D24').Value = '10.85'
$ws.Range('E24').Value = '  -1.37%  '
$ws.Range('E25').Value = '  -1.91%  '
$ws.Range('D26').Value = '157.37'
$ws.Range('E26').Value = '  -1.23%  '
$ws.Range('E27').Value = '  -3.10%  '
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('D29').Value = '5.319'
$ws.Range('E29').Value = '  -3.42%  '
$ws.Range('D30').Value = '118.96'
$ws.Range('E30').Value = '  -2.23%  '
$ws.Range('D31').Value = '0.9531'
$ws.Range('E31').Value = '  -2.63%  '
$ws.Range('D32').Value = '0.09298'
$ws.Range('E32').Value = '  -2.79%  '
$ws.Range('D33').Value = '3.565'
$ws.Range('E33').Value = '  -1.83%  '
$ws.Range('D34').Value = '5.241'
$ws.Range('E34').Value = '  -1.99%  '
$ws.Range('D35').Value = '1.318'
$ws.Range('E35').Value = '  -3.34%  '
$ws.Range('D36').Value = '0.05936'
$ws.Range('E36').Value = '  -2.20%  '
$ws.Range('D37').Value = '0.02190'
$ws.Range('E37').Value = '  -2.78%  '
$ws.Range('D38').Value = '8.071'
$ws.Range('E38').Value = '  -1.58%  '
$ws.Range('D39').Value = '1.157'
$ws.Range('E39').Value = '  -4.08%  '
$ws.Range('D40').Value = '0.5793'
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('D41').Value = '0.1841'
$ws.Range('E41').Value = '  -2.54%  '
$ws.Range('E42').Value = '  -3.05%  '
$ws.Range('D43').Value = '1.260'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '12.05'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.5484'
$ws.Range('E45').Value = '  -3.12%  '
$ws.Range('D46').Value = '1.868'
$ws.Range('E46').Value = '  -3.38%  '
$ws.Range('D47').Value = '0.06659'
$ws.Range('E47').Value = '  -2.25%  '
$ws.Range('E48').Value = '  -2.12%  '
$ws.Range('D49').Value = '1.042'
$ws.Range('E49').Value = '  -2.85%  '
$ws.Range('D50').Value = '1.002'
$ws.Range('E50').Value = '  -0.98%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.00000000288'
$ws.Range('E51').Value = '  -0.11%  '

$ws.Range('D2:D51').Style = 'Normal'
